$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Brunnen (... risiegen Brunnen reichen)" -> restructure:
#       "Brunnen (... bis zu " + "riesigen" + " Brunnen reichen)" + "wa"(proofErr)
# ---------------------------------------------------------------------------

# 1a) Change the (still proofErr-wrapped) "risiegen" run's text to "wa" -- a plain
#     word-for-word replace keeps the surrounding spellStart/spellEnd proofErr tags.
$d.Content.Find.Execute("risiegen", $false, $false, $false, $false, $false, $true, 1, $false, "wa", 2) | Out-Null

# 1b) Delete the old " Brunnen reichen)" run (it will be reinserted earlier below).
$delRange = $d.Content
$delRange.Find.Execute(" Brunnen reichen)") | Out-Null
$delRange.Text = ""

# 1c) Insert the two new runs ("riesigen" and " Brunnen reichen)") right after the
#     unchanged prefix run, before the proofErr-wrapped "wa" run. A temporary,
#     distinct font size is applied right after each insertion so the new text
#     does not get silently merged into a neighbouring run of identical
#     formatting; the size is restored to the original (18 half-points = 9pt)
#     once both runs exist as independent runs.
$prefixRange = $d.Content
$prefixRange.Find.Execute("Brunnen (Kann von einem kleinen Trinkbrunnen bis zu ") | Out-Null
$insPos = $prefixRange.End

$ins1 = $d.Range($insPos, $insPos)
$ins1.InsertBefore("riesigen")
$run1 = $d.Range($insPos, $d.Content.End)
$run1.Find.Execute("riesigen") | Out-Null
$run1.Font.Size = 77
$run1Start = $run1.Start
$run1End = $run1.End

$ins2 = $d.Range($run1End, $run1End)
$ins2.InsertBefore(" Brunnen reichen)")
$run2 = $d.Range($run1End, $d.Content.End)
$run2.Find.Execute(" Brunnen reichen)") | Out-Null
$run2.Font.Size = 66
$run2Start = $run2.Start
$run2End = $run2.End

$d.Range($run1Start, $run1End).Font.Size = 9
$d.Range($run2Start, $run2End).Font.Size = 9

# ---------------------------------------------------------------------------
# 2) "Abfalleimer (Um anzugeben ... hinzufügen)" -> "Abfalleimer"
# ---------------------------------------------------------------------------
$oldEimer = "Abfalleimer (Um anzugeben welche Art von Abfall hier deponiert wird, kann man den Tag „waste=art“ hinzufügen)"
$d.Content.Find.Execute($oldEimer, $false, $false, $false, $false, $false, $true, 1, $false, "Abfalleimer", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Abfallcontainer (Um anzugeben ... hinzufügen)" -> "Abfallcontainer"
#    followed by the relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$oldContainer = "Abfallcontainer (Um anzugeben welche Art von Abfall hier deponiert wird, kann man den Tag „waste=art“ hinzufügen)"
$d.Content.Find.Execute($oldContainer, $false, $false, $false, $false, $false, $true, 1, $false, "Abfallcontainer", 2) | Out-Null

$containerRange = $d.Content
$containerRange.Find.Execute("Abfallcontainer") | Out-Null
$bmPos = $containerRange.End

# The original "_GoBack" bookmark (right after "Erfassungskatalog" near the start
# of the document) is implicitly removed by Bookmarks.Add below: Word allows only
# one bookmark per name, so re-adding "_GoBack" moves it to the new location.
$bmIns = $d.Range($bmPos, $bmPos)
$bmIns.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $bmIns) | Out-Null
$bmRange = $d.Bookmarks("_GoBack").Range
$bmRange.Text = ""
